# Apply a new "Office Theme" colour scheme to the presentation's design
# (this mirrors picking a different theme from the PowerPoint Design
# gallery: ppt/theme/theme1.xml - the theme used by the slide master /
# all slides - gets the stock "Office" palette instead of the previous
# "Integral" / Red Violet palette).
#
# The 12 theme colours, in the standard OOXML <a:clrScheme> order
# (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink), expressed as COM
# RGB() long values (0x00BBGGRR) so they round-trip through
# ColorFormat.RGB the same way real PowerPoint COM automation would
# write them:
#   dk1      #000000 ->         0
#   lt1      #FFFFFF ->  16777215
#   dk2      #44546A ->   6968388
#   lt2      #E7E6E6 ->  15132391
#   accent1  #5B9BD5 ->  13998939
#   accent2  #ED7D31 ->   3243501
#   accent3  #A5A5A5 ->  10855845
#   accent4  #FFC000 ->     49407
#   accent5  #4472C4 ->  12874308
#   accent6  #70AD47 ->   4697456
#   hlink    #0563C1 ->  12673797
#   folHlink #954F72 ->   7491477

$p = $ppt.ActivePresentation

$design = $p.Designs.Item(1)
$slideMaster = $design.SlideMaster
$colorScheme = $slideMaster.Theme.ThemeColorScheme

$newThemeColors = @(0, 16777215, 6968388, 15132391, 13998939, 3243501, 10855845, 49407, 12874308, 4697456, 12673797, 7491477)

for ($i = 1; $i -le $colorScheme.Count; $i++) {
    $colorScheme.Item($i).RGB = $newThemeColors[$i - 1]
}
